$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 17:03"

# Country stats table (rows 4-218), sorted descending by "Casos totales" (col B).
# Refreshed data caused some countries to change rank, so both the figures and the
# row order are rewritten here to match the refreshed/re-sorted dataset.
$data = @(
    ,("Estados Unidos", 1191965, 3843, 178671, 944579, 16139, 117, 68715)
    ,("España", 248301, 1179, 151633, 71240, 2254, 164, 25428)
    ,("Italia", 210717, 0, 81654, 100179, 1501, 0, 28884)
    ,("Reino Unido", 186599, 0, 0, 157809, 1559, 0, 28446)
    ,("Francia", 168693, 0, 50784, 93014, 3819, 0, 24895)
    ,("Alemania", 165745, 81, 132700, 26179, 1979, 0, 6866)
    ,("Rusia", 145268, 10581, 18095, 125817, 2300, 76, 1356)
    ,("Turquia", 126045, 0, 63151, 59497, 1424, 0, 3397)
    ,("Brasil", 101826, 679, 42991, 51784, 8318, 26, 7051)
    ,("Iran", 98647, 1223, 79379, 12991, 2676, 74, 6277)
    ,("China", 82880, 3, 77766, 481, 33, 0, 4633)
    ,("Canada", 59474, 0, 24908, 30884, 557, 0, 3682)
    ,("Belgica", 50267, 361, 12378, 29965, 655, 80, 7924)
    ,("Peru", 45928, 0, 13550, 31092, 679, 0, 1286)
    ,("India", 42836, 331, 11782, 29659, 0, 4, 1395)
    ,("Paises Bajos", 40770, 199, 0, 35438, 683, 26, 5082)
    ,("Suiza", 29981, 76, 24500, 3719, 141, 0, 1762)
    ,("Ecuador", 29538, 0, 3300, 24674, 149, 0, 1564)
    ,("Arabia Saudita", 28656, 1645, 4476, 23989, 139, 7, 191)
    ,("Portugal", 25524, 242, 1712, 22749, 143, 20, 1063)
    ,("Mexico", 23471, 1383, 13447, 7870, 378, 93, 2154)
    ,("Suecia", 22721, 404, 1005, 18947, 455, 90, 2769)
    ,("Irlanda", 21506, 0, 13386, 6817, 98, 0, 1303)
    ,("Pakistan", 20884, 800, 5635, 14773, 111, 19, 476)
    ,("Chile", 19663, 0, 10041, 9362, 408, 0, 260)
    ,("Singapur", 18778, 573, 1457, 17303, 22, 0, 18)
    ,("Bielorrusia", 17489, 784, 3259, 14127, 92, 4, 103)
    ,("Israel", 16237, 29, 9858, 6145, 93, 2, 234)
    ,("Catar", 16191, 640, 1810, 14369, 72, 0, 12)
    ,("Austria", 15621, 24, 13316, 1705, 111, 2, 600)
    ,("Japon", 14877, 0, 3981, 10409, 321, 0, 487)
    ,("Emiratos Arabes Unidos", 14163, 0, 2763, 11274, 1, 0, 126)
    ,("Polonia", 13937, 244, 4095, 9159, 160, 5, 683)
    ,("Rumania", 13512, 349, 5269, 7440, 243, 13, 803)
    ,("Ucrania", 12331, 418, 1619, 10409, 160, 15, 303)
    ,("Indonesia", 11587, 395, 1954, 8769, 0, 19, 864)
    ,("Corea del Sur", 10801, 8, 9217, 1332, 55, 2, 252)
    ,("Banglades", 10143, 688, 1209, 8752, 1, 5, 182)
    ,("Dinamarca", 9670, 147, 7088, 2089, 57, 9, 493)
    ,("Serbia", 9557, 93, 1574, 7786, 53, 4, 197)
    ,("Filipinas", 9485, 262, 1315, 7547, 31, 16, 623)
    ,("Republica Dominicana", 8235, 281, 1771, 6118, 144, 13, 346)
    ,("Noruega", 7884, 37, 32, 7638, 27, 3, 214)
    ,("Chequia", 7781, 0, 3592, 3940, 58, 1, 249)
    ,("Colombia", 7668, 0, 1722, 5606, 122, 0, 340)
    ,("Panama", 7197, 107, 641, 6356, 91, 3, 200)
    ,("Australia", 6825, 24, 5859, 871, 28, 0, 95)
    ,("Sudafrica", 6783, 0, 2549, 4103, 36, 0, 131)
    ,("Egipto", 6465, 0, 1562, 4474, 0, 0, 429)
    ,("Malasia", 6353, 55, 4484, 1764, 28, 0, 105)
    ,("Finlandia", 5327, 73, 3500, 1587, 49, 10, 240)
    ,("Kuwait", 5278, 295, 1947, 3291, 79, 2, 40)
    ,("Marruecos", 5000, 97, 1565, 3258, 1, 3, 177)
    ,("Argentina", 4783, 0, 1354, 3183, 157, 0, 246)
    ,("Argelia", 4474, 0, 1936, 2075, 22, 0, 463)
    ,("Moldavia", 4248, 127, 1423, 2697, 237, 3, 128)
    ,("Kazajistan", 4028, 108, 1134, 2867, 40, 0, 27)
    ,("Luxemburgo", 3824, 0, 3379, 349, 22, 0, 96)
    ,("Barein", 3383, 0, 1718, 1657, 1, 0, 8)
    ,("Hungria", 3035, 37, 630, 2054, 55, 11, 351)
    ,("Tailandia", 2987, 18, 2740, 193, 61, 0, 54)
    ,("Afganistan", 2894, 190, 397, 2407, 7, 5, 90)
    ,("Oman", 2637, 69, 816, 1809, 17, 0, 12)
    ,("Grecia", 2626, 0, 1374, 1108, 37, 0, 144)
    ,("Nigeria", 2558, 0, 400, 2071, 4, 0, 87)
    ,("Armenia", 2507, 121, 1071, 1397, 10, 4, 39)
    ,("Irak", 2296, 0, 1490, 709, 0, 0, 97)
    ,("Uzbekistan", 2181, 32, 1370, 801, 8, 0, 10)
    ,("Ghana", 2169, 0, 229, 1922, 4, 0, 18)
    ,("Croacia", 2101, 5, 1522, 499, 15, 1, 80)
    ,("Camerun", 2077, 0, 953, 1060, 12, 0, 64)
    ,("Azerbaiyan", 1932, 0, 1441, 466, 17, 0, 25)
    ,("Bosnia y Herzegovina", 1926, 69, 855, 993, 4, 1, 78)
    ,("Islandia", 1799, 0, 1717, 72, 0, 0, 10)
    ,("Estonia", 1703, 3, 259, 1389, 6, 0, 55)
    ,("Bulgaria", 1652, 34, 321, 1253, 40, 5, 78)
    ,("Cuba", 1649, 0, 827, 755, 13, 0, 67)
    ,("Bolivia", 1594, 124, 166, 1352, 3, 5, 76)
    ,("Guinea", 1586, 0, 405, 1174, 0, 0, 7)
    ,("Republica de Macedonia", 1518, 7, 992, 441, 21, 1, 85)
    ,("Nueva Zelanda", 1487, 0, 1276, 191, 0, 0, 20)
    ,("Eslovenia", 1439, 0, 241, 1101, 20, 1, 97)
    ,("Lituania", 1419, 9, 638, 735, 17, 0, 46)
    ,("Eslovaquia", 1413, 5, 643, 745, 7, 1, 25)
    ,("Costa de Marfil", 1398, 0, 653, 728, 0, 0, 17)
    ,("Senegal", 1271, 89, 415, 846, 4, 1, 10)
    ,("Republica de Yibuti", 1116, 4, 713, 401, 0, 0, 2)
    ,("Honduras", 1055, 45, 118, 855, 10, 6, 82)
    ,("Hong Kong", 1041, 1, 900, 137, 1, 0, 4)
    ,("Tunez", 1013, 0, 328, 643, 25, 0, 42)
    ,("Letonia", 896, 17, 348, 532, 4, 0, 16)
    ,("Republica de Chipre", 874, 2, 296, 563, 15, 0, 15)
    ,("Kirguistan", 830, 35, 575, 245, 13, 0, 10)
    ,("Albania", 803, 8, 543, 229, 7, 0, 31)
    ,("Niger", 750, 0, 518, 196, 0, 0, 36)
    ,("Principado de Andorra", 748, 0, 493, 210, 17, 0, 45)
    ,("Libano", 740, 3, 200, 515, 43, 0, 25)
    ,("Costa Rica", 739, 0, 386, 347, 6, 0, 6)
    ,("Somalia", 722, 0, 44, 646, 2, 0, 32)
    ,("Sri Lanka", 721, 3, 194, 519, 1, 1, 8)
    ,("Crucero", 712, 0, 645, 54, 4, 0, 13)
    ,("Guatemala", 703, 15, 72, 614, 5, 0, 17)
    ,("Consejo Danes para los Refugiados", 682, 8, 80, 568, 0, 1, 34)
    ,("Sudan", 678, 86, 61, 576, 0, 0, 41)
    ,("Burkina Faso", 662, 0, 540, 77, 0, 0, 45)
    ,("Uruguay", 655, 0, 442, 196, 10, 0, 17)
    ,("Mayotte", 650, 0, 235, 409, 4, 0, 6)
    ,("Georgia", 593, 4, 223, 361, 6, 0, 9)
    ,("San Marino", 582, 0, 86, 455, 5, 0, 41)
    ,("Mali", 580, 17, 223, 328, 0, 2, 29)
    ,("El Salvador", 555, 65, 180, 363, 3, 1, 12)
    ,("Maldivas", 527, 0, 18, 508, 2, 0, 1)
    ,("Kenia", 490, 25, 173, 293, 2, 0, 24)
    ,("Tanzania", 480, 0, 167, 297, 7, 0, 16)
    ,("Malta", 480, 3, 399, 77, 1, 0, 4)
    ,("Jamaica", 469, 6, 38, 422, 2, 1, 9)
    ,("Jordania", 461, 0, 367, 85, 5, 0, 9)
    ,("Taiwan", 438, 6, 334, 98, 0, 0, 6)
    ,("Reunion", 423, 0, 300, 123, 2, 0, 0)
    ,("Paraguay", 396, 26, 126, 260, 7, 0, 10)
    ,("Venezuela", 357, 0, 158, 189, 1, 0, 10)
    ,("Estado de Palestina", 353, 0, 77, 274, 0, 0, 2)
    ,("Gabon", 335, 0, 85, 245, 1, 0, 5)
    ,("Mauricio", 332, 0, 316, 6, 3, 0, 10)
    ,("Montenegro", 323, 1, 249, 66, 2, 0, 8)
    ,("Isla de Man", 321, 0, 271, 28, 21, 0, 22)
    ,("Guinea Ecuatorial", 315, 0, 13, 301, 0, 0, 1)
    ,("Vietnam", 271, 0, 219, 52, 8, 0, 0)
    ,("Ruanda", 259, 0, 124, 135, 0, 0, 0)
    ,("Guinea-Bisau", 257, 0, 19, 237, 0, 0, 1)
    ,("Tayikistan", 230, 102, 0, 227, 0, 1, 3)
    ,("Congo", 229, 0, 25, 195, 0, 0, 9)
    ,("Islas Feroe", 187, 0, 185, 2, 0, 0, 0)
    ,("Martinica", 179, 0, 83, 82, 5, 0, 14)
    ,("Sierra Leona", 178, 12, 37, 132, 0, 1, 9)
    ,("Cabo Verde", 175, 10, 37, 136, 0, 0, 2)
    ,("Liberia", 166, 8, 58, 90, 0, 0, 18)
    ,("Birmania", 161, 6, 49, 106, 0, 0, 6)
    ,("Guadalupe", 152, 0, 95, 45, 6, 0, 12)
    ,("Madagascar", 149, 0, 99, 50, 1, 0, 0)
    ,("Gibraltar", 144, 0, 132, 12, 0, 0, 0)
    ,("Etiopia", 140, 5, 75, 62, 0, 0, 3)
    ,("Brunei", 138, 0, 130, 7, 2, 0, 1)
    ,("Zambia", 137, 13, 78, 56, 1, 0, 3)
    ,("Guayana Francesa", 128, 0, 98, 29, 2, 0, 1)
    ,("Togo", 126, 2, 74, 43, 0, 0, 9)
    ,("Camboya", 122, 0, 120, 2, 1, 0, 0)
    ,("Republica del Chad", 117, 0, 39, 68, 0, 0, 10)
    ,("Trinidad yTobago", 116, 0, 96, 12, 0, 0, 8)
    ,("Bermudas", 115, 0, 51, 57, 4, 0, 7)
    ,("Suazilandia", 112, 0, 12, 99, 0, 0, 1)
    ,("Aruba", 100, 0, 81, 17, 4, 0, 2)
    ,("Benin", 96, 6, 50, 44, 0, 0, 2)
    ,("Monaco", 95, 0, 78, 13, 1, 0, 4)
    ,("Uganda", 89, 0, 55, 34, 0, 0, 0)
    ,("Haiti", 88, 3, 10, 69, 0, 1, 9)
    ,("Republica de Africa Central", 85, 13, 10, 75, 0, 0, 0)
    ,("Bahamas", 83, 0, 24, 48, 1, 0, 11)
    ,("Guyana", 82, 0, 22, 51, 2, 0, 9)
    ,("Barbados", 82, 1, 44, 31, 4, 0, 7)
    ,("Liechtenstein", 82, 0, 55, 26, 0, 0, 1)
    ,("Mozambique", 80, 0, 19, 61, 0, 0, 0)
    ,("San Martin (Parte Holandesa)", 76, 0, 44, 19, 7, 0, 13)
    ,("Nepal", 75, 0, 16, 59, 0, 0, 0)
    ,("Islas Caimanes", 74, 0, 10, 63, 3, 0, 1)
    ,("Libia", 63, 0, 22, 38, 0, 0, 3)
    ,("Polinesia Francesa", 58, 0, 51, 7, 1, 0, 0)
    ,("Sudan del Sur", 46, 0, 0, 46, 0, 0, 0)
    ,("Macao", 45, 0, 39, 6, 1, 0, 0)
    ,("Siria", 44, 0, 27, 14, 0, 0, 3)
    ,("Malaui", 41, 2, 9, 29, 1, 0, 3)
    ,("Mongolia", 40, 1, 12, 28, 0, 0, 0)
    ,("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
    ,("Eritrea", 39, 0, 26, 13, 0, 0, 0)
    ,("San Martin (Parte Francesa)", 38, 0, 27, 8, 3, 0, 3)
    ,("Angola", 35, 0, 11, 22, 0, 0, 2)
    ,("Zimbabue", 34, 0, 5, 25, 0, 0, 4)
    ,("Guam", 32, 0, 0, 31, 0, 0, 1)
    ,("Antigua y Barbuda", 25, 0, 15, 7, 1, 0, 3)
    ,("Timor Oriental", 24, 0, 20, 4, 0, 0, 0)
    ,("Botsuana", 23, 0, 8, 14, 0, 0, 1)
    ,("Granada", 21, 0, 13, 8, 4, 0, 0)
    ,("Laos", 19, 0, 9, 10, 0, 0, 0)
    ,("Fiyi", 18, 0, 14, 4, 0, 0, 0)
    ,("Santa Lucia", 18, 0, 15, 3, 0, 0, 0)
    ,("Belice", 18, 0, 13, 3, 1, 0, 2)
    ,("Nueva Caledonia", 18, 0, 17, 1, 1, 0, 0)
    ,("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
    ,("San Vicente y las Granadinas", 17, 1, 9, 8, 0, 0, 0)
    ,("Gambia", 17, 0, 9, 7, 0, 0, 1)
    ,("Santo Tome y Principe", 16, 0, 4, 11, 0, 0, 1)
    ,("Namibia", 16, 0, 8, 8, 0, 0, 0)
    ,("Dominica", 16, 0, 13, 3, 0, 0, 0)
    ,("Curazao", 16, 0, 13, 2, 0, 0, 1)
    ,("San Cristobal y Nieves", 15, 0, 8, 7, 0, 0, 0)
    ,("Burundi", 15, 0, 7, 7, 0, 0, 1)
    ,("Nicaragua", 15, 1, 7, 3, 0, 2, 5)
    ,("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
    ,("Islas Turcas y Caicos", 12, 0, 6, 5, 0, 0, 1)
    ,("Santa Sede", 11, 0, 2, 9, 0, 0, 0)
    ,("Seychelles", 11, 0, 6, 5, 0, 0, 0)
    ,("Montserrat", 11, 0, 7, 3, 1, 0, 1)
    ,("Groenlandia", 11, 0, 11, 0, 0, 0, 0)
    ,("Yemen", 10, 0, 1, 7, 0, 0, 2)
    ,("Surinam", 10, 0, 9, 0, 0, 0, 1)
    ,("Papua Nueva Guinea", 8, 0, 6, 2, 0, 0, 0)
    ,("Mauritania", 8, 0, 6, 1, 0, 0, 1)
    ,("Butan", 7, 0, 5, 2, 0, 0, 0)
    ,("Bonaire, San Eustaquio y Saba", 6, 0, 0, 6, 0, 0, 0)
    ,("Islas Virgenes Britanicas", 6, 0, 3, 2, 0, 0, 1)
    ,("Sahara Occidental", 6, 0, 5, 1, 0, 0, 0)
    ,("San Bartolome", 6, 0, 6, 0, 0, 0, 0)
    ,("Comoras", 3, 0, 0, 3, 0, 0, 0)
    ,("Anguila", 3, 0, 3, 0, 0, 0, 0)
    ,("San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

"Updated " + ($r - 4) + " rows"
